# Logout Methode im UserController erstellt und Layout für Anmeldung angepasst
$wb = $excel.ActiveWorkbook

$wsUserStories = $wb.Worksheets.Item("UserStories")
$wsTasks       = $wb.Worksheets.Item("Tasks")

# --- Tasks sheet: "Logout Methode im UserController erstellen" row (row 22) ---
# Mark it as "done" (was referencing the stray "b" shared string) and stamp
# the completion date, matching the sibling rows above it.
$wsTasks.Range("D21").Copy($wsTasks.Range("D22"))
$wsTasks.Range("C22").Value = "done"
$wsTasks.Range("D22").Value = 43530

# --- View / selection state ---
# Anmeldung (login) work is now the focus: select the Tasks sheet's new
# current entry row, then switch over to / select on UserStories.
$wsTasks.Activate() | Out-Null
$wsTasks.Range("D24").Select() | Out-Null

$wsUserStories.Activate() | Out-Null
$wsUserStories.Range("B14").Select() | Out-Null
